$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume number, report date range) ---
$ws.Range("A8").Value = "Volume 32   Number  22"
$ws.Range("C9").Value = "Report Covering the Week  5/26/2025  Through  6/1/2025"

# --- Numeric cell updates ---
$ws.Range("C15").Value = 2
$ws.Range("D15").Value = 2
$ws.Range("F15").Value = 6
$ws.Range("G15").Value = 6
$ws.Range("H15").Value = 0
$ws.Range("I15").Value = 17
$ws.Range("J15").Value = 18
$ws.Range("K15").Value = -5.555555555555
$ws.Range("L15").Value = 21.428571428571
$ws.Range("M15").Value = 183.333333333333
$ws.Range("N15").Value = 21.428571428571
$ws.Range("C16").Value = 5
$ws.Range("D16").Value = 2
$ws.Range("E16").Value = 150
$ws.Range("F16").Value = 29
$ws.Range("G16").Value = 13
$ws.Range("H16").Value = 123.076923076923
$ws.Range("I16").Value = 94
$ws.Range("J16").Value = 106
$ws.Range("K16").Value = -11.320754716981
$ws.Range("L16").Value = -25.396825396825
$ws.Range("M16").Value = -6.930693069306
$ws.Range("N16").Value = -71.515151515151
$ws.Range("C17").Value = 6
$ws.Range("D17").Value = 13
$ws.Range("E17").Value = -53.846153846153
$ws.Range("G17").Value = 44
$ws.Range("H17").Value = -31.818181818181
$ws.Range("I17").Value = 179
$ws.Range("J17").Value = 187
$ws.Range("K17").Value = -4.278074866310
$ws.Range("L17").Value = 9.146341463414
$ws.Range("M17").Value = 10.493827160493
$ws.Range("N17").Value = -1.104972375690
$ws.Range("C18").Value = 7
$ws.Range("D18").Value = 5
$ws.Range("E18").Value = 40
$ws.Range("F18").Value = 20
$ws.Range("G18").Value = 19
$ws.Range("H18").Value = 5.263157894736
$ws.Range("I18").Value = 87
$ws.Range("J18").Value = 82
$ws.Range("K18").Value = 6.097560975609
$ws.Range("L18").Value = 7.407407407407
$ws.Range("M18").Value = -9.375
$ws.Range("N18").Value = -77.225130890052
$ws.Range("C19").Value = 7
$ws.Range("D19").Value = 9
$ws.Range("E19").Value = -22.222222222222
$ws.Range("F19").Value = 30
$ws.Range("H19").Value = -11.764705882352
$ws.Range("I19").Value = 136
$ws.Range("J19").Value = 171
$ws.Range("K19").Value = -20.467836257309
$ws.Range("L19").Value = -16.564417177914
$ws.Range("M19").Value = 27.102803738317
$ws.Range("N19").Value = -16.564417177914
$ws.Range("C20").Value = 3
$ws.Range("D20").Value = 7
$ws.Range("E20").Value = -57.142857142857
$ws.Range("F20").Value = 19
$ws.Range("G20").Value = 18
$ws.Range("H20").Value = 5.555555555555
$ws.Range("I20").Value = 79
$ws.Range("J20").Value = 80
$ws.Range("K20").Value = -1.25
$ws.Range("L20").Value = -44.755244755244
$ws.Range("M20").Value = 132.352941176471
$ws.Range("N20").Value = -52.121212121212
$ws.Range("D21").Value = 38
$ws.Range("E21").Value = -21.052631578947
$ws.Range("F21").Value = 134
$ws.Range("G21").Value = 134
$ws.Range("H21").Value = 0
$ws.Range("I21").Value = 592
$ws.Range("J21").Value = 646
$ws.Range("K21").Value = -8.359133126934
$ws.Range("L21").Value = -14.697406340057
$ws.Range("M21").Value = 16.99604743083
$ws.Range("N21").Value = -52.564102564102
$ws.Range("L22").Value = -20
$ws.Range("M22").Value = -27.272727272727
$ws.Range("D23").Value = 3
$ws.Range("G23").Value = 5
$ws.Range("J23").Value = 15
$ws.Range("K23").Value = -86.666666666666
$ws.Range("L23").Value = -81.818181818181
$ws.Range("C24").Value = 16
$ws.Range("D24").Value = 14
$ws.Range("E24").Value = 14.285714285714
$ws.Range("F24").Value = 60
$ws.Range("G24").Value = 79
$ws.Range("H24").Value = -24.050632911392
$ws.Range("I24").Value = 358
$ws.Range("J24").Value = 330
$ws.Range("K24").Value = 8.484848484848
$ws.Range("L24").Value = 12.578616352201
$ws.Range("M24").Value = 56.331877729257
$ws.Range("C25").Value = 3
$ws.Range("D25").Value = 6
$ws.Range("E25").Value = -50
$ws.Range("F25").Value = 15
$ws.Range("G25").Value = 23
$ws.Range("H25").Value = -34.782608695652
$ws.Range("I25").Value = 95
$ws.Range("J25").Value = 90
$ws.Range("K25").Value = 5.555555555555
$ws.Range("L25").Value = -7.766990291262
$ws.Range("C26").Value = 5
$ws.Range("D26").Value = 11
$ws.Range("E26").Value = -54.545454545454
$ws.Range("G26").Value = 57
$ws.Range("H26").Value = -22.807017543859
$ws.Range("I26").Value = 252
$ws.Range("J26").Value = 269
$ws.Range("K26").Value = -6.319702602230
$ws.Range("L26").Value = 26
$ws.Range("M26").Value = 4.564315352697
$ws.Range("C27").Value = 2
$ws.Range("D27").Value = 2
$ws.Range("F27").Value = 7
$ws.Range("G27").Value = 7
$ws.Range("H27").Value = 0
$ws.Range("I27").Value = 28
$ws.Range("J27").Value = 32
$ws.Range("K27").Value = -12.5
$ws.Range("L27").Value = 16.666666666666
$ws.Range("D28").Value = 5
$ws.Range("E28").Value = -100
$ws.Range("F28").Value = 5
$ws.Range("G28").Value = 10
$ws.Range("H28").Value = -50
$ws.Range("J28").Value = 49
$ws.Range("K28").Value = -20.408163265306
$ws.Range("E29").Value = -100
$ws.Range("G29").Value = 2
$ws.Range("H29").Value = 0
$ws.Range("J29").Value = 10
$ws.Range("K29").Value = -80
$ws.Range("L29").Value = -80
$ws.Range("M29").Value = -75
$ws.Range("N29").Value = -93.548387096774
$ws.Range("E30").Value = -100
$ws.Range("G30").Value = 2
$ws.Range("H30").Value = 0
$ws.Range("J30").Value = 9
$ws.Range("K30").Value = -77.777777777777
$ws.Range("L30").Value = -77.777777777777
$ws.Range("M30").Value = -75
$ws.Range("N30").Value = -93.333333333333

# --- Cells changing from numeric to text placeholder "0" (style copied from C31) ---
$ws.Range("C31").Copy()
$ws.Range("C22").NumberFormat = "@"
$ws.Range("C22").Value = "0"
$ws.Range("C22").PasteSpecial(-4122)
$ws.Range("C28").NumberFormat = "@"
$ws.Range("C28").Value = "0"
$ws.Range("C28").PasteSpecial(-4122)
$ws.Range("C29").NumberFormat = "@"
$ws.Range("C29").Value = "0"
$ws.Range("C29").PasteSpecial(-4122)
$ws.Range("C30").NumberFormat = "@"
$ws.Range("C30").Value = "0"
$ws.Range("C30").PasteSpecial(-4122)
